$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the dates between the "8-abril" rows (2,3,4) and "9-abril" rows (5,6,7)
$ws.Range("D2").Value = 44295
$ws.Range("D3").Value = 44295
$ws.Range("D4").Value = 44295

$ws.Range("D5").Value = 44294
$ws.Range("D6").Value = 44294
$ws.Range("D7").Value = 44294

# Swap the Volumen (M) values between row 3 (Primera) and row 6 (Primera)
$ws.Range("M3").Value = 200
$ws.Range("M6").Value = 240
